$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44371
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 1800
$ws.Range("O2").Value = 1800
$ws.Range("P2").Value = 1800
$ws.Range("S2").Value = 1800

# Row 3
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 1200
$ws.Range("O3").Value = 1200
$ws.Range("P3").Value = 1200
$ws.Range("S3").Value = 1200

# Row 4
$ws.Range("D4").Value = 44336
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 1500
$ws.Range("O4").Value = 1500
$ws.Range("P4").Value = 1500
$ws.Range("S4").Value = 1500

# Row 5
$ws.Range("D5").Value = 44309
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 1600
$ws.Range("O5").Value = 1600
$ws.Range("P5").Value = 1600
$ws.Range("S5").Value = 1600

# Row 6
$ws.Range("D6").Value = 44292
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 14000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("S6").Value = 1400
$ws.Range("T6").Value = 10

# Row 7
$ws.Range("D7").Value = 44400
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 1500
$ws.Range("O7").Value = 1500
$ws.Range("P7").Value = 1500
$ws.Range("Q7").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44195
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("S8").Value = 1500
$ws.Range("T8").Value = 10

# Row 11
$ws.Range("D11").Value = 44343
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 1700
$ws.Range("O11").Value = 1700
$ws.Range("P11").Value = 1700
$ws.Range("Q11").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S11").Value = 1700
$ws.Range("T11").Value = 1
